$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52 held the customer whose phone number was stored as text
# "09876543" (leading zero). Re-type that cell as a plain number,
# matching how every other row on the sheet stores the phone column.
# Birthday (B52, blank) and points (C52, already 0) are left untouched.
$ws.Cells.Item(52, 1).Value = 9876543

# Append a new row 53 for that same customer, preserving the original
# text-formatted phone number (with its leading zero) and (re)setting
# their points to 0.00.
$ws.Cells.Item(53, 1).NumberFormat = "@"
$ws.Cells.Item(53, 1).Value = "09876543"
$ws.Cells.Item(53, 1).Style = "Normal"
$ws.Cells.Item(53, 2).Value = ""
$ws.Cells.Item(53, 3).Value = 0
